$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.2792131854384032, 0.9817265900740069, 3.212386962184873, 1.792313299115105, 1.788047957535663, 51)
    3  = @(0.04913906157728941, 0.7455142430313828, 1.431191262711348, 1.196324062581435, 1.20823713588917, 47)
    4  = @(0.09986602815088967, 0.8843907116147273, 2.216810884689522, 1.488895860928333, 1.501603274128746, 47)
    5  = @(0.07441293746091567, 0.9042953237502608, 2.143707891555164, 1.464140666587457, 1.478771582361615, 45)
    6  = @(-0.08417250974636116, 0.7926005748427621, 1.965095792784742, 1.401818744625974, 1.415466677756357, 44)
    7  = @(-0.05842981899582228, 0.890948235336648, 1.810325433388112, 1.345483345637586, 1.365722831991239, 32)
    8  = @(-0.05865351751528118, 0.8879676780668504, 1.798954600371419, 1.341251132477218, 1.36211777908385, 31)
    9  = @(0.09549959853857792, 0.7676328100216852, 1.068282700646155, 1.033577621974351, 1.068005883826769, 14)
    10 = @(0.2725186903762485, 0.7176745655965183, 1.025465364696565, 1.01265263772755, 1.042633439540573, 8)
    11 = @(0.1847752594325322, 0.7699318286940575, 1.072099526920768, 1.035422390583074, 1.176411283195583, 4)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
